# Applies the "added account creatn details" edit:
#  - renames Sheet2 -> accountCreationdetails
#  - fills in account-creation-details data on that sheet
#  - adds a mailto hyperlink on B2
#  - sets a few column widths
#  - makes the new sheet the active / selected tab

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

$ws2.Name = "accountCreationdetails"

# --- Row 1 headers, and Row 2 data, written in the same order the
#     original author entered them (this drives shared-string order) ---

# Batch 1: headers for the "core" columns (A,B,C,D, I,J,K, M,N,O,P)
$ws2.Range("A1").Value = "Name"
$ws2.Range("B1").Value = "Email"
$ws2.Range("C1").Value = "Gender"
$ws2.Range("D1").Value = "Namee"
$ws2.Range("E1").Value = "Password"
$ws2.Range("I1").Value = "Firstname"
$ws2.Range("J1").Value = "LastName"
$ws2.Range("K1").Value = "Address"
$ws2.Range("M1").Value = "State"
$ws2.Range("N1").Value = "City "
$ws2.Range("O1").Value = "Zipcode"
$ws2.Range("P1").Value = "Number"

# Batch 2: row 2 data for those same columns
$ws2.Range("A2").Value = "raja"
$ws2.Range("B2").Value = "BabbLe@gmail.com"
$ws2.Range("C2").Value = "mr"
$ws2.Range("D2").Value = "raja"
$ws2.Range("E2").Value = "ran1231"
$ws2.Range("F2").Value = 15
$ws2.Range("I2").Value = "raja"
$ws2.Range("J2").Value = "kaja"
$ws2.Range("K2").Value = "hyderabad"
$ws2.Range("M2").Value = "matagat"
$ws2.Range("N2").Value = "njhggd"
$ws2.Range("O2").Value = 637465
$ws2.Range("P2").Value = 7344253663

# Batch 3: the remaining headers (Day / Month / Year / Country), added later
$ws2.Range("F1").Value = "Day"
$ws2.Range("G1").Value = "Month"
$ws2.Range("H1").Value = "Year"
$ws2.Range("L1").Value = "Country"

# Batch 4: Country value
$ws2.Range("L2").Value = "India"

# Numeric year
$ws2.Range("H2").Value = 2005

# Batch 5: Month value (entered last of all new strings)
$ws2.Range("G2").Value = "August"

# --- Hyperlink + style on B2 ---
$ws2.Hyperlinks.Add($ws2.Range("B2"), "mailto:BabbLe@gmail.com") | Out-Null
$ws2.Range("B2").Style = "Hyperlink"

# --- Column widths ---
$ws2.Columns.Item(2).ColumnWidth = 17.333333333333336
$ws2.Columns.Item(8).ColumnWidth = 9.833333333333332
$ws2.Columns.Item(9).ColumnWidth = 9.333333333333332
$ws2.Columns.Item(12).ColumnWidth = 10.333333333333332
$ws2.Columns.Item(15).ColumnWidth = 10.0

# --- Selection on the new sheet, then make it the active tab ---
$ws2.Range("M1:P2").Select() | Out-Null
$ws2.Activate() | Out-Null

Write-Output "done"
